$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving exact formatting
# (e.g. leading/trailing zeros, multi-dot separators) without Excel
# auto-converting it to a number, and without leaving a stray
# NumberFormat style applied to the cell afterwards.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 2 4 "63.991.13"
Set-TextValue 2 5 "  -1.31%  "
# Row 3
Set-TextValue 3 4 "3.398.67"
Set-TextValue 3 5 "  -1.55%  "
# Row 4
Set-TextValue 4 5 "  +0.03%  "
# Row 5
Set-TextValue 5 4 "571.32"
Set-TextValue 5 5 "  -0.50%  "
# Row 6
Set-TextValue 6 4 "162.79"
Set-TextValue 6 5 "  +2.05%  "
# Row 7
Set-TextValue 7 5 "  +0.04%  "
# Row 8
Set-TextValue 8 4 "3.400.50"
Set-TextValue 8 5 "  -1.41%  "
# Row 9
Set-TextValue 9 4 "0.548"
Set-TextValue 9 5 "  -5.10%  "
# Row 10
Set-TextValue 10 4 "7.29"
Set-TextValue 10 5 "  +1.44%  "
# Row 11
Set-TextValue 11 5 "  -2.40%  "
# Row 12
Set-TextValue 12 4 "0.420"
Set-TextValue 12 5 "  -4.65%  "
# Row 13
Set-TextValue 13 4 "3.985.97"
Set-TextValue 13 5 "  -1.47%  "
# Row 14
Set-TextValue 14 5 "  +0.61%  "
# Row 15
Set-TextValue 15 4 "26.92"
Set-TextValue 15 5 "  -2.74%  "
# Row 16
Set-TextValue 16 4 "0.0000172"
Set-TextValue 16 5 "  -2.99%  "
# Row 17
Set-TextValue 17 4 "64.029.99"
Set-TextValue 17 5 "  -1.33%  "
# Row 18
Set-TextValue 18 4 "3.398.71"
Set-TextValue 18 5 "  -0.90%  "
# Row 19
Set-TextValue 19 4 "6.11"
Set-TextValue 19 5 "  -1.77%  "
# Row 20
Set-TextValue 20 4 "13.58"
Set-TextValue 20 5 "  -1.34%  "
# Row 21
Set-TextValue 21 4 "376.43"
Set-TextValue 21 5 "  -0.50%  "
# Row 22
Set-TextValue 22 4 "7.77"
Set-TextValue 22 5 "  -2.51%  "
# Row 23
Set-TextValue 23 5 "  +0.05%  "
# Row 24
Set-TextValue 24 4 "70.28"
Set-TextValue 24 5 "  -2.68%  "
# Row 25
Set-TextValue 25 4 "0.511"
Set-TextValue 25 5 "  -4.57%  "
# Row 26
Set-TextValue 26 4 "0.0000114"
Set-TextValue 26 5 "  -6.25%  "
# Row 27
Set-TextValue 27 5 "  -4.25%  "
# Row 28
Set-TextValue 28 4 "0.178"
Set-TextValue 28 5 "  -0.46%  "
# Row 29
Set-TextValue 29 4 "0.999"
Set-TextValue 29 5 "  -0.06%  "
# Row 30
Set-TextValue 30 4 "6.10"
Set-TextValue 30 5 "  +0.50%  "
# Row 31
Set-TextValue 31 4 "1.39"
Set-TextValue 31 5 "  -4.50%  "
# Row 32
Set-TextValue 32 4 "2.00"
Set-TextValue 32 5 "  -0.62%  "
# Row 33
Set-TextValue 33 4 "22.76"
Set-TextValue 33 5 "  -1.98%  "
# Row 34
Set-TextValue 34 4 "7.02"
Set-TextValue 34 5 "  +0.34%  "
# Row 35
Set-TextValue 35 4 "1.48"
Set-TextValue 35 5 "  -5.67%  "
# Row 36
Set-TextValue 36 4 "159.80"
Set-TextValue 36 5 "  -0.65%  "
# Row 37
Set-TextValue 37 5 "  +9.25%  "
# Row 38
Set-TextValue 38 4 "1.80"
Set-TextValue 38 5 "  -4.11%  "
# Row 39
Set-TextValue 39 4 "0.0720"
Set-TextValue 39 5 "  -3.92%  "
# Row 40
Set-TextValue 40 2 "EnergySwap"
Set-TextValue 40 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 40 4 "25.72"
Set-TextValue 40 5 "  -1.98%  "
# Row 41
Set-TextValue 41 2 "OKB"
Set-TextValue 41 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 41 4 "42.83"
Set-TextValue 41 5 "  -0.36%  "
# Row 42
Set-TextValue 42 2 "Maker"
Set-TextValue 42 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 42 4 "2.741.68"
Set-TextValue 42 5 "  -5.42%  "
# Row 43
Set-TextValue 43 4 "26.23"
Set-TextValue 43 5 "  +0.38%  "
# Row 44
Set-TextValue 44 4 "6.41"
Set-TextValue 44 5 "  -2.40%  "
# Row 45
Set-TextValue 45 4 "4.36"
Set-TextValue 45 5 "  -3.77%  "
# Row 46
Set-TextValue 46 4 "0.0305"
Set-TextValue 46 5 "  -2.31%  "
# Row 47
Set-TextValue 47 4 "2.40"
Set-TextValue 47 5 "  +0.87%  "
# Row 48
Set-TextValue 48 4 "327.38"
Set-TextValue 48 5 "  +2.07%  "
# Row 49
Set-TextValue 49 5 "  -4.35%  "
# Row 50
Set-TextValue 50 5 "  -2.01%  "
# Row 51
Set-TextValue 51 4 "6.25"
Set-TextValue 51 5 "  -3.43%  "
